$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 4 (current "Layout" row), shifting
# everything below down by two rows.
$ws.Rows("4:5").Insert()

# Populate the new rows in column B.
$ws.Range("B4").Value = "Add media player on page"
$ws.Range("B5").Value = "Add links to YouTube pages"

# Update the active selection to match the target state.
$ws.Range("B6").Select()
